$d = $word.ActiveDocument

# 1) Mint the "Tabelraster" (Table Grid) table style used by the new table.
#    (Word normally creates this built-in style the first time a table is
#    formatted with the "Table Grid" gallery style.)
$tableGridStyle = $d.Styles.Add("Tabelraster", 3)
$tableGridStyle.NameLocal = "Table Grid"
$tableGridStyle.BaseStyle = "Standaardtabel"
$tableGridStyle.Priority = 39
$tableGridStyle.ParagraphFormat.SpaceAfter = 0
$tableGridStyle.ParagraphFormat.LineSpacingRule = 0

# 2) Locate the existing "GDPR" Kop1 heading paragraph and expand the range
#    to cover the whole paragraph (including its paragraph mark) so that
#    InsertXML below replaces it wholesale with the new block of content:
#    "Non-functional requirements" heading + table + page break + the
#    (re-created) "GDPR" heading.
$headingRange = $d.Content
$found = $headingRange.Find.Execute("GDPR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingRange.Expand(4) | Out-Null

$newBlockXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Kop1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Non-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>functional</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>requirements</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="Tabelraster"/><w:tblW w:w="0" w:type="auto"/><w:tblInd w:w="0" w:type="dxa"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="805"/><w:gridCol w:w="2340"/><w:gridCol w:w="5917"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="805" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>NF</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>01</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>Schaalbaarheid</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5917" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>E</w:t></w:r><w:r><w:t xml:space="preserve">r </w:t></w:r><w:r><w:t xml:space="preserve">wordt </w:t></w:r><w:r><w:t xml:space="preserve">rekening gehouden met de mogelijkheid om de applicatie makkelijk uit te kunnen breiden. </w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="805" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>NF</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>02</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>Schaalbaarheid</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5917" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>De</w:t></w:r><w:r><w:t xml:space="preserve"> applicatie kan automatisch op en afschalen wanneer dit nodig is.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="805" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>NF</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>Performance</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5917" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>De snelheid wordt grotendeels gewaarborgd door een juiste architectuur. De applicatie dient te blijven werken onder zware load.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="805" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>NF</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>Robuustheid</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5917" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t xml:space="preserve">De applicatie gaat juist om met </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>errors</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Hij vangt deze af en zorgt ervoor dat het </w:t></w:r><w:r><w:t>niet zichtbaar</w:t></w:r><w:r><w:t xml:space="preserve"> is voor de gebruiker.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="805" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>NF</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>Gedistribueerde data</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5917" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>In het design en in de architectuur moet rekening worden gehouden met het juist opslaan van data</w:t></w:r><w:r><w:t xml:space="preserve"> zie ook GDPR.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="805" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>NF</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>6</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>Security</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5917" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>E</w:t></w:r><w:r><w:t>r</w:t></w:r><w:r><w:t xml:space="preserve"> wordt</w:t></w:r><w:r><w:t xml:space="preserve"> bij het design en de architectuur rekening gehouden met </w:t></w:r><w:r><w:t>de</w:t></w:r><w:r><w:t xml:space="preserve"> risico</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> van OWASP</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">Ook worden er aan </w:t></w:r><w:r><w:t>basisprincipes</w:t></w:r><w:r><w:t xml:space="preserve"> gehouden van security.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="805" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>NF</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>7</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t>Privacy</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5917" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:r><w:t xml:space="preserve">Er wordt rekening gehouden met de </w:t></w:r><w:r><w:t>privacyregels</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>over, van en voor</w:t></w:r><w:r><w:t xml:space="preserve"> het opslaan van data.</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Kop1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>GDPR</w:t></w:r></w:p></w:body></w:wordDocument>
'@

$headingRange.InsertXML($newBlockXml)

# 3) Tidy up the "De volgende data ... niet opgeslagen:" paragraph: collapse
#    the three runs (incl. the gramStart/gramEnd proofErr-wrapped "wordt")
#    into a single plain run with the same text.
$gdprIntro = $d.Content
$gdprIntro.Find.Execute("De volgende data", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$gdprIntro.Expand(4) | Out-Null
$gdprIntro.Text = "De volgende data wordt in verband met GDPR niet opgeslagen:"
